$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 546.7765096666667
$ws.Range("H2").Value = 1640.329529
$ws.Range("I2").Value = 0.6285526459909564
$ws.Range("J2").Value = 0.6285526459909564
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 91918.504629986
$ws.Range("R2").Value = 827266.541669874
$ws.Range("S2").Value = 0.1875715654405607
$ws.Range("T2").Value = 0.1875715654405608
$ws.Range("G3").Value = 546.7765096666667
$ws.Range("H3").Value = 1640.329529
$ws.Range("I3").Value = 0.6285526459909564
$ws.Range("J3").Value = 0.6285526459909564
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 89127.98059172212
$ws.Range("R3").Value = 802151.8253254991
$ws.Range("S3").Value = 0.1818771411854699
$ws.Range("T3").Value = 0.1818771411854699
$ws.Range("G4").Value = 546.7765096666667
$ws.Range("H4").Value = 1640.329529
$ws.Range("I4").Value = 0.6285526459909564
$ws.Range("J4").Value = 0.6285526459909564
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 90761.36296064913
$ws.Range("R4").Value = 816852.2666458421
$ws.Range("S4").Value = 0.1852102686023698
$ws.Range("T4").Value = 0.1852102686023698
$ws.Range("G5").Value = 546.7765096666667
$ws.Range("H5").Value = 1640.329529
$ws.Range("I5").Value = 0.6285526459909564
$ws.Range("J5").Value = 0.6285526459909564
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 36211.22264540161
$ws.Range("R5").Value = 325901.0038086144
$ws.Range("S5").Value = 0.07389367076255597
$ws.Range("T5").Value = 0.07389367076255597
$ws.Range("I6").Value = 0.1861770314550556
$ws.Range("J6").Value = 0.1861770314550556
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 27226.22271491449
$ws.Range("R6").Value = 245036.0044342304
$ws.Range("S6").Value = 0.05555861941213393
$ws.Range("T6").Value = 0.05555861941213394
$ws.Range("I7").Value = 0.1861770314550556
$ws.Range("J7").Value = 0.1861770314550556
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.05387193332399084
$ws.Range("T7").Value = 0.05387193332399084
$ws.Range("I8").Value = 0.1861770314550556
$ws.Range("J8").Value = 0.1861770314550556
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 26883.47783531822
$ws.Range("R8").Value = 241951.300517864
$ws.Range("S8").Value = 0.05485920427400259
$ws.Range("T8").Value = 0.05485920427400259
$ws.Range("I9").Value = 0.1861770314550556
$ws.Range("J9").Value = 0.1861770314550556
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 10725.74903705991
$ws.Range("R9").Value = 96531.74133353916
$ws.Range("S9").Value = 0.02188727444492827
$ws.Range("T9").Value = 0.02188727444492827
$ws.Range("G10").Value = 160.630483
$ws.Range("H10").Value = 481.891449
$ws.Range("I10").Value = 0.1846544489960017
$ws.Range("J10").Value = 0.1846544489960017
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 27003.56276159993
$ws.Range("R10").Value = 243032.0648543994
$ws.Range("S10").Value = 0.0551042530560639
$ws.Range("T10").Value = 0.05510425305606392
$ws.Range("G11").Value = 160.630483
$ws.Range("H11").Value = 481.891449
$ws.Range("I11").Value = 0.1846544489960017
$ws.Range("J11").Value = 0.1846544489960017
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 26183.77036714849
$ws.Range("R11").Value = 235653.9333043364
$ws.Range("S11").Value = 0.05343136092860257
$ws.Range("T11").Value = 0.05343136092860257
$ws.Range("G12").Value = 160.630483
$ws.Range("H12").Value = 481.891449
$ws.Range("I12").Value = 0.1846544489960017
$ws.Range("J12").Value = 0.1846544489960017
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 26663.62089877499
$ws.Range("R12").Value = 239972.5880889749
$ws.Range("S12").Value = 0.05441055783522093
$ws.Range("T12").Value = 0.05441055783522093
$ws.Range("G13").Value = 160.630483
$ws.Range("H13").Value = 481.891449
$ws.Range("I13").Value = 0.1846544489960017
$ws.Range("J13").Value = 0.1846544489960017
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 10638.03232347602
$ws.Range("R13").Value = 95742.29091128414
$ws.Range("S13").Value = 0.0217082771761143
$ws.Range("T13").Value = 0.0217082771761143
$ws.Range("G14").Value = 0.5357470000000001
$ws.Range("H14").Value = 1.607241
$ws.Range("I14").Value = 0.0006158735579862568
$ws.Range("J14").Value = 0.0006158735579862568
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 90.06433566435135
$ws.Range("R14").Value = 810.579020979162
$ws.Range("S14").Value = 0.0001837878944933949
$ws.Range("T14").Value = 0.000183787894493395
$ws.Range("G15").Value = 0.5357470000000001
$ws.Range("H15").Value = 1.607241
$ws.Range("I15").Value = 0.0006158735579862568
$ws.Range("J15").Value = 0.0006158735579862568
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 87.33010173970968
$ws.Range("R15").Value = 785.9709156573871
$ws.Range("S15").Value = 0.0001782083374761193
$ws.Range("T15").Value = 0.0001782083374761193
$ws.Range("G16").Value = 0.5357470000000001
$ws.Range("H16").Value = 1.607241
$ws.Range("I16").Value = 0.0006158735579862568
$ws.Range("J16").Value = 0.0006158735579862568
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 88.93053571691001
$ws.Range("R16").Value = 800.3748214521901
$ws.Range("S16").Value = 0.0001814742294496251
$ws.Range("T16").Value = 0.0001814742294496251
$ws.Range("G17").Value = 0.5357470000000001
$ws.Range("H17").Value = 1.607241
$ws.Range("I17").Value = 0.0006158735579862568
$ws.Range("J17").Value = 0.0006158735579862568
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 35.48077423888034
$ws.Range("R17").Value = 319.3269681499231
$ws.Range("S17").Value = 0.00007240309656711741
$ws.Range("T17").Value = 0.00007240309656711741
